$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text value to a cell while preventing Excel from
# auto-converting plain decimal-looking strings (e.g. "236.35") into
# numbers, and without leaving any lasting NumberFormat/style change.
function Set-TextValue($cell, $value) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = $origStyle
}

$ws.Range("D2").Value = '30.314.10'
$ws.Range("E2").Value = '  -0.02%  '

$ws.Range("D3").Value = '1.869.22'
$ws.Range("E3").Value = '  +0.14%  '

$ws.Range("E4").Value = '  -0.10%  '

Set-TextValue $ws.Range("D5") '236.35'
$ws.Range("E5").Value = '  +0.31%  '

Set-TextValue $ws.Range("D7") '0.4710'
$ws.Range("E7").Value = '  +0.71%  '

Set-TextValue $ws.Range("D8") '0.2902'
$ws.Range("E8").Value = '  +2.23%  '

Set-TextValue $ws.Range("D9") '0.06626'
$ws.Range("E9").Value = '  +1.61%  '

Set-TextValue $ws.Range("D10") '21.71'
$ws.Range("E10").Value = '  -0.18%  '

Set-TextValue $ws.Range("D11") '0.08019'
$ws.Range("E11").Value = '  +1.03%  '

Set-TextValue $ws.Range("D12") '97.31'
$ws.Range("E12").Value = '  -0.02%  '

$ws.Range("D13").Value = '1.877.30'
$ws.Range("E13").Value = '  +0.52%  '

Set-TextValue $ws.Range("D14") '5.154'
$ws.Range("E14").Value = '  +0.14%  '

$ws.Range("E15").Value = '  +1.72%  '

Set-TextValue $ws.Range("D16") '273.94'
$ws.Range("E16").Value = '  -1.88%  '

$ws.Range("D17").Value = '30.318.99'
$ws.Range("E17").Value = '  +0.03%  '

Set-TextValue $ws.Range("D18") '14.09'
$ws.Range("E18").Value = '  +5.66%  '

Set-TextValue $ws.Range("D19") '0.000007740'
$ws.Range("E19").Value = '  +5.96%  '

$ws.Range("E20").Value = '  -0.10%  '

$ws.Range("D21").Value = '2.120.88'
$ws.Range("E21").Value = '  +0.22%  '

$ws.Range("E22").Value = '  -1.68%  '

$ws.Range("E23").Value = '  -0.13%  '

Set-TextValue $ws.Range("D24") '6.219'
$ws.Range("E24").Value = '  +1.07%  '

Set-TextValue $ws.Range("D25") '167.75'
$ws.Range("E25").Value = '  +0.81%  '

Set-TextValue $ws.Range("D26") '9.286'
$ws.Range("E26").Value = '  +1.30%  '

Set-TextValue $ws.Range("D27") '18.98'
$ws.Range("E27").Value = '  -0.60%  '

$ws.Range("E28").Value = '  +1.69%  '

$ws.Range("E29").Value = '  -1.21%  '

$ws.Range("E30").Value = '  +2.35%  '

Set-TextValue $ws.Range("D31") '4.366'
$ws.Range("E31").Value = '  -0.77%  '

$ws.Range("E32").Value = '  -0.89%  '

Set-TextValue $ws.Range("D33") '4.094'
$ws.Range("E33").Value = '  +0.37%  '

$ws.Range("E34").Value = '  -0.74%  '

$ws.Range("E35").Value = '  +0.28%  '

Set-TextValue $ws.Range("D36") '0.7035'
$ws.Range("E36").Value = '  -0.46%  '

Set-TextValue $ws.Range("D37") '2.702'
$ws.Range("E37").Value = '  -0.40%  '

Set-TextValue $ws.Range("D38") '0.01884'
$ws.Range("E38").Value = '  +1.10%  '

Set-TextValue $ws.Range("D39") '2.646'
$ws.Range("E39").Value = '  +2.60%  '

Set-TextValue $ws.Range("D40") '6.327'
$ws.Range("E40").Value = '  +0.10%  '

Set-TextValue $ws.Range("D41") '73.27'
$ws.Range("E41").Value = '  -1.77%  '

Set-TextValue $ws.Range("D42") '1.964'
$ws.Range("E42").Value = '  -0.13%  '

Set-TextValue $ws.Range("D43") '0.8419'
$ws.Range("E43").Value = '  -0.99%  '

Set-TextValue $ws.Range("D44") '0.4171'
$ws.Range("E44").Value = '  -0.40%  '

Set-TextValue $ws.Range("D45") '0.9998'
$ws.Range("E45").Value = '  -0.15%  '

Set-TextValue $ws.Range("D46") '103.64'
$ws.Range("E46").Value = '  +0.32%  '

Set-TextValue $ws.Range("D47") '7.143'
$ws.Range("E47").Value = '  -0.71%  '

Set-TextValue $ws.Range("D48") '9.227'
$ws.Range("E48").Value = '  -1.11%  '

Set-TextValue $ws.Range("D49") '934.66'
$ws.Range("E49").Value = '  -3.16%  '

Set-TextValue $ws.Range("D51") '0.05668'
$ws.Range("E51").Value = '  +0.35%  '
